# lsh_coding.xlsx — "Added new unit categores and fixed bug."
#
# 1) Append 4 new unit-category rows to the "lsh_unit_categories" sheet
#    (rows 28-31), mirroring the existing Bráðamóttaka/Göngudeild pattern.
# 2) Make "lsh_unit_categories" the active sheet/tab (fixing the stray
#    tabSelected that had been left on "lsh_text_out_categories").

$wb = $excel.ActiveWorkbook

$wsUnit = $wb.Worksheets.Item("lsh_unit_categories")

$wsUnit.Range("A28").Value = "Bráðaöldrunarlækningadeild (Fv-B4)                                                                                           "
$wsUnit.Range("B28").Value = "Bráðamóttaka"
$wsUnit.Range("C28").Value = "emergency_room"
$wsUnit.Range("D28").Value = "home"
$wsUnit.Range("E28").Value = 1

$wsUnit.Range("A29").Value = "Fv-A3 GD Svefnrannsókna"
$wsUnit.Range("B29").Value = "Göngudeild"
$wsUnit.Range("C29").Value = "outpatient_clinic"
$wsUnit.Range("D29").Value = "home"
$wsUnit.Range("E29").Value = 1

$wsUnit.Range("A30").Value = "Hb-21B GD Fósturgreiningadeild"
$wsUnit.Range("B30").Value = "Göngudeild"
$wsUnit.Range("C30").Value = "outpatient_clinic"
$wsUnit.Range("D30").Value = "home"
$wsUnit.Range("E30").Value = 1

$wsUnit.Range("A31").Value = "Hb-31E GD Geðsviðs"
$wsUnit.Range("B31").Value = "Göngudeild"
$wsUnit.Range("C31").Value = "outpatient_clinic"
$wsUnit.Range("D31").Value = "home"
$wsUnit.Range("E31").Value = 1

# Make this sheet the active one (this also clears tabSelected on whatever
# sheet used to be selected, e.g. lsh_text_out_categories) and land the
# selection/scroll position where the author left it.
$wsUnit.Activate()
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
$wsUnit.Range("E32").Select()
